$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly data row at row 996 (shifts the existing rows
# 996-1048 down to 997-1049, matching the published dimension A1:R1049).
$ws.Rows.Item(996).Insert()

$ws.Cells.Item(996, 1).Value = 8
$ws.Cells.Item(996, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(996, 3).Value = "Coquimbo"
$ws.Cells.Item(996, 4).Value = 45267
$ws.Cells.Item(996, 5).Value = 4
$ws.Cells.Item(996, 6).Value = 100112043
$ws.Cells.Item(996, 7).Value = "Pepino ensalada"
$ws.Cells.Item(996, 8).Value = "Sin especificar"
$ws.Cells.Item(996, 9).Value = "Primera"
$ws.Cells.Item(996, 10).Value = 520
$ws.Cells.Item(996, 11).Value = 14000
$ws.Cells.Item(996, 12).Value = 15000
$ws.Cells.Item(996, 13).Value = 14500
$ws.Cells.Item(996, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(996, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(996, 16).Value = 242
$ws.Cells.Item(996, 17).Value = 60
$ws.Cells.Item(996, 18).Value = "Hortaliza"
